$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("All Published Values")
$summary = $wb.Worksheets.Item("Daily Summary")

# Append new row 17 to "All Published Values" with a new BOC USD rate
# publication. Force text formatting before writing so values that look
# like dates/numbers ("2026-01-02", "697.85", ...) are stored as literal
# text, matching every other row in the sheet, then reset the style back
# to Normal so no unwanted formatting/number-format sticks to the cells.
$ws.Range("A17:J17").NumberFormat = "@"
$ws.Range("A17").Value = "2026-01-02"
$ws.Range("B17").Value = "2026-01-02 22:31:02"
$ws.Range("C17").Value = "697.85"
$ws.Range("D17").Value = "697.85"
$ws.Range("E17").Value = "700.79"
$ws.Range("F17").Value = "700.79"
$ws.Range("G17").Value = "702.88"
$ws.Range("H17").Value = "2026/01/02 22:31:02"
$ws.Range("I17").Value = "2026-01-02 14:40:29"
$ws.Range("J17").Value = "https://www.bankofchina.com/sourcedb/whpj/enindex_1619.html"
$ws.Range("A17:J17").Style = "Normal"

# Grow the autofilter to cover the newly added row. The sheet already has
# an autofilter, so first switch it off, then re-apply it over the full
# A1:J17 range (calling .AutoFilter() on a range that already carries a
# filter just toggles it off, so clear first).
$ws.AutoFilterMode = $false
$ws.Range("A1:J17").AutoFilter() | Out-Null

# Keep the hidden _FilterDatabase defined name for this sheet in sync with
# the resized autofilter range.
$wb.Names.Item("All Published Values!_FilterDatabase").RefersTo = "='All Published Values'!`$A`$1:`$J`$17"

# Daily Summary: bump the publishes count for 2026-01-02 from 15 to 16
# now that an extra publication was captured.
$summary.Range("B4").Value = 16
